# Lightning talk deck — "Update some literature review"
#
# Slide 6 ("Python Implementation"), Content Placeholder 2: revise the
# literature-review bullets — reword the first two points and add a new
# bullet about real-time data support.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(6)

# Find the "Content Placeholder 2" shape by name (robust to shape ordering).
$contentShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.Name -eq "Content Placeholder 2") {
        $contentShape = $sh
    }
}

$tr = $contentShape.TextFrame.TextRange

# 1) Reword: "Mainstream algorithm does not seem to support time-dependent
#    routing" -> "Mainstream algorithm is a black box to me"
$hit = $tr.Find("Mainstream algorithm does not seem to support time-dependent routing", 0)
$hit.Text = "Mainstream algorithm is a black box to me"

# 2) Reword: "If they do, please tell me, e.g. open trip planner" ->
#    "If they do support time-dependent routing, please tell me, e.g. open
#    trip planner"
$hit2 = $tr.Find("If they do, please tell me, e.g. open trip planner", 0)
$hit2.Text = "If they do support time-dependent routing, please tell me, e.g. open trip planner"

# 3) Add a new bullet right after "Much more flexible for scientific
#    calculation": "Support for achieved real-time data"
# (Use the Paragraphs() accessor rather than Find() here: InsertAfter on a
# paragraph range creates a genuine new <a:p>, whereas InsertAfter on a
# Find() sub-run range merges the text into the same run/paragraph.)
$paraCount = ($tr.Text -split "`r").Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.Trim() -eq "Much more flexible for scientific calculation") {
        $para.InsertAfter("`rSupport for achieved real-time data") | Out-Null
        break
    }
}
